# Generate Report for Handback
# For the "08f4a62e-5dac-4afe-87b3-811ce7cafee7" handback entry (row 5 on both the
# zh-cn and de-de sheets), fill in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" / "Error Detail" columns, which the handback report
# generator had left blank, and widen a few columns to match the other wide,
# free-text columns.

$wb = $excel.ActiveWorkbook

$currentUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/839bae1745b299482fac4f1190bdee217237ff58/e2e/08f4a62e-5dac-4afe-87b3-811ce7cafee7.md"
$latestUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1fa1dd52a7f3fc5ac6038c759c3107815eaa84b1/e2e/08f4a62e-5dac-4afe-87b3-811ce7cafee7.md"
$errorDetail = "The version of handback file is not the latest, current: $currentUrl, latest: $latestUrl."
$targetMdName = "08f4a62e-5dac-4afe-87b3-811ce7cafee7.md"

function Update-HandbackRow {
    param(
        $ws,
        [string]$handbackXlf,
        [string]$handbackDatetime
    )

    # I5: "Latest Target File" - hyperlink to the current (non-latest) commit of the handback md file
    $ws.Hyperlinks.Add($ws.Range("I5"), $currentUrl, "", "", $targetMdName) | Out-Null

    # J5: "Latest Handback File"
    $ws.Range("J5").Value = $handbackXlf

    # K5: "Latest Handback DateTime"
    $ws.Range("K5").Value = $handbackDatetime

    # P5: "Error Detail"
    $ws.Range("P5").Value = $errorDetail

    # Widen the columns that now hold longer free-text content, matching the
    # width (40) already used for the other wide text columns (A, G).
    $ws.Columns.Item(9).ColumnWidth = 39.17
    $ws.Columns.Item(10).ColumnWidth = 39.17
    $ws.Columns.Item(16).ColumnWidth = 39.17
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HandbackRow -ws $wsZhCn `
    -handbackXlf "08f4a62e-5dac-4afe-87b3-811ce7cafee7.5ae7c0831bea55190f7d7c1547df7d4493662f2d.zh-cn.xlf" `
    -handbackDatetime "2016-10-25 02:10:47"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HandbackRow -ws $wsDeDe `
    -handbackXlf "08f4a62e-5dac-4afe-87b3-811ce7cafee7.5ae7c0831bea55190f7d7c1547df7d4493662f2d.de-de.xlf" `
    -handbackDatetime "2016-10-25 02:11:05"
